$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows (11-13) - Resolving-Mac as a sending cluster is removed
$ws.Rows("11:13").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf22"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1094766666666667
$ws.Range("H2").Value = 0.32843
$ws.Range("I2").Value = 0.4072979538991744
$ws.Range("J2").Value = 0.4072979538991744
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8155003333333334
$ws.Range("N2").Value = 2.446501
$ws.Range("O2").Value = 0.1910612426590028
$ws.Range("P2").Value = 0.1910612426590029
$ws.Range("Q2").Value = 0.08927825815888889
$ws.Range("R2").Value = 0.80350432343
$ws.Range("S2").Value = 0.07781885320444552
$ws.Range("T2").Value = 0.07781885320444552

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf22"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1094766666666667
$ws.Range("H3").Value = 0.32843
$ws.Range("I3").Value = 0.4072979538991744
$ws.Range("J3").Value = 0.4072979538991744
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.333134333333334
$ws.Range("N3").Value = 9.999403000000001
$ws.Range("O3").Value = 0.7809105179307759
$ws.Range("P3").Value = 0.780910517930776
$ws.Range("Q3").Value = 0.3649004363655556
$ws.Range("R3").Value = 3.28410392729
$ws.Range("S3").Value = 0.3180632561315496
$ws.Range("T3").Value = 0.3180632561315496

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf22"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1094766666666667
$ws.Range("H4").Value = 0.32843
$ws.Range("I4").Value = 0.4072979538991744
$ws.Range("J4").Value = 0.4072979538991744
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.119632
$ws.Range("N4").Value = 0.358896
$ws.Range("O4").Value = 0.02802823941022116
$ws.Range("P4").Value = 0.02802823941022117
$ws.Range("Q4").Value = 0.01309691258666667
$ws.Range("R4").Value = 0.11787221328
$ws.Range("S4").Value = 0.01141584456317928
$ws.Range("T4").Value = 0.01141584456317928

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Fgf22"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09509666666666666
$ws.Range("H5").Value = 0.28529
$ws.Range("I5").Value = 0.3537984753764744
$ws.Range("J5").Value = 0.3537984753764744
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8155003333333334
$ws.Range("N5").Value = 2.446501
$ws.Range("O5").Value = 0.1910612426590028
$ws.Range("P5").Value = 0.1910612426590029
$ws.Range("Q5").Value = 0.07755136336555556
$ws.Range("R5").Value = 0.69796227029
$ws.Range("S5").Value = 0.06759717635628981
$ws.Range("T5").Value = 0.06759717635628981

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Fgf22"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.09509666666666666
$ws.Range("H6").Value = 0.28529
$ws.Range("I6").Value = 0.3537984753764744
$ws.Range("J6").Value = 0.3537984753764744
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.333134333333334
$ws.Range("N6").Value = 9.999403000000001
$ws.Range("O6").Value = 0.7809105179307759
$ws.Range("P6").Value = 0.780910517930776
$ws.Range("Q6").Value = 0.3169699646522222
$ws.Range("R6").Value = 2.85272968187
$ws.Range("S6").Value = 0.2762849506493615
$ws.Range("T6").Value = 0.2762849506493615

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Fgf22"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.09509666666666666
$ws.Range("H7").Value = 0.28529
$ws.Range("I7").Value = 0.3537984753764744
$ws.Range("J7").Value = 0.3537984753764744
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.119632
$ws.Range("N7").Value = 0.358896
$ws.Range("O7").Value = 0.02802823941022116
$ws.Range("P7").Value = 0.02802823941022117
$ws.Range("Q7").Value = 0.01137660442666667
$ws.Range("R7").Value = 0.10238943984
$ws.Range("S7").Value = 0.00991634837082306
$ws.Range("T7").Value = 0.009916348370823062

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf22"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.06421433333333333
$ws.Range("H8").Value = 0.192643
$ws.Range("I8").Value = 0.2389035707243512
$ws.Range("J8").Value = 0.2389035707243512
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8155003333333334
$ws.Range("N8").Value = 2.446501
$ws.Range("O8").Value = 0.1910612426590028
$ws.Range("P8").Value = 0.1910612426590029
$ws.Range("Q8").Value = 0.05236681023811111
$ws.Range("R8").Value = 0.471301292143
$ws.Range("S8").Value = 0.04564521309826751
$ws.Range("T8").Value = 0.04564521309826752

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf22"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.06421433333333333
$ws.Range("H9").Value = 0.192643
$ws.Range("I9").Value = 0.2389035707243512
$ws.Range("J9").Value = 0.2389035707243512
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.333134333333334
$ws.Range("N9").Value = 9.999403000000001
$ws.Range("O9").Value = 0.7809105179307759
$ws.Range("P9").Value = 0.780910517930776
$ws.Range("Q9").Value = 0.2140349991254445
$ws.Range("R9").Value = 1.926314992129
$ws.Range("S9").Value = 0.1865623111498649
$ws.Range("T9").Value = 0.1865623111498649

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf22"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.06421433333333333
$ws.Range("H10").Value = 0.192643
$ws.Range("I10").Value = 0.2389035707243512
$ws.Range("J10").Value = 0.2389035707243512
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.119632
$ws.Range("N10").Value = 0.358896
$ws.Range("O10").Value = 0.02802823941022116
$ws.Range("P10").Value = 0.02802823941022117
$ws.Range("Q10").Value = 0.007682089125333334
$ws.Range("R10").Value = 0.069138802128
$ws.Range("S10").Value = 0.006696046476218819
$ws.Range("T10").Value = 0.00669604647621882
